$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target contents for A2:B12 after the edit (row order was reshuffled,
# and a few label values were updated along the way).
$data = @(
  @("", "Ica Informa"),
  @("Se declara admisible 08 días. Concede ONI.", "Admisibles"),
  @("Agréguese Extraordinariamente a la tabla del 05 de marzo", "Agréguese A Tabla"),
  @("Rechazada sin costas", "Rechazada"),
  @("Ev. Informe. En relación", "Evacua Informe"),
  @("Inadmisible", "Inadmisible/Omite"),
  @("Se Pronuncia Incompetencia", "Incompetencia"),
  @("Dese Cuenta Admisibilidad", "Dese Cuenta"),
  @("Se declara admisible, 10 días. Concede ONI.", "Admisibles"),
  @("Concede ampliación de plazo por el término de 08 días hábiles", "Ica Amplia Plazo"),
  @("Atendiendo al tiempo transcurrido y que la recurrida no ha evacuado el informe dese cuenta para los fines que diere lugar", "Evacua Informe")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
